$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.419217109680176
$ws.Range("B1").Value = 2.793301582336426
$ws.Range("C1").Value = 3.000523805618286
$ws.Range("D1").Value = 3.604333400726318
$ws.Range("E1").Value = 1.650992870330811
